$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.913.88"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "'3.482.88"
$ws.Range("E3").Value = "  -2.49%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'601.60"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("D6").Value = "'148.25"
$ws.Range("E6").Value = "  -4.71%  "

$ws.Range("D7").Value = "'3.480.74"
$ws.Range("E7").Value = "  -2.51%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("E11").Value = "  +3.64%  "

$ws.Range("D12").Value = "'0.422"
$ws.Range("E12").Value = "  -3.77%  "

$ws.Range("D13").Value = "'0.0000212"
$ws.Range("E13").Value = "  -4.25%  "

$ws.Range("D14").Value = "'4.072.93"
$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("E15").Value = "  -6.35%  "

$ws.Range("D16").Value = "'3.495.41"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("D17").Value = "'66.957.13"
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "'6.39"
$ws.Range("E19").Value = "  -5.21%  "

$ws.Range("D20").Value = "'10.14"
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("E21").Value = "  -5.07%  "

$ws.Range("D22").Value = "'434.18"
$ws.Range("E22").Value = "  -4.68%  "

$ws.Range("E23").Value = "  -6.09%  "

$ws.Range("D24").Value = "'79.19"
$ws.Range("E24").Value = "  +0.82%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").Value = "'3.620.76"
$ws.Range("E26").Value = "  -2.45%  "

$ws.Range("E27").Value = "  -10.30%  "

$ws.Range("D28").Value = "'9.81"
$ws.Range("E28").Value = "  -7.19%  "

$ws.Range("E29").Value = "  -9.79%  "

$ws.Range("D30").Value = "'2.48"
$ws.Range("E30").Value = "  -3.51%  "

$ws.Range("E31").Value = "  -6.67%  "

$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "'25.35"
$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("D35").Value = "'3.473.48"
$ws.Range("E35").Value = "  -2.55%  "

$ws.Range("D36").Value = "'5.92"
$ws.Range("E36").Value = "  -7.47%  "

$ws.Range("E37").Value = "  -6.67%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").Value = "'7.89"
$ws.Range("E39").Value = "  -4.45%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("D41").Value = "'173.30"
$ws.Range("E41").Value = "  -4.55%  "

$ws.Range("D42").Value = "'0.0884"
$ws.Range("E42").Value = "  -3.82%  "

$ws.Range("E43").Value = "  -12.90%  "

$ws.Range("E44").Value = "  -3.85%  "

$ws.Range("D45").Value = "'0.899"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").Value = "'29.09"
$ws.Range("E46").Value = "  -6.85%  "

$ws.Range("D47").Value = "'46.45"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("E48").Value = "  -7.20%  "

$ws.Range("D49").Value = "'7.44"
$ws.Range("E49").Value = "  -4.74%  "

$ws.Range("D50").Value = "'2.41"
$ws.Range("E50").Value = "  -10.24%  "

$ws.Range("D51").Value = "'0.970"
$ws.Range("E51").Value = "  -5.32%  "
